$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the section labels down one row, into the first data row of each block.
$ws.Range("A4").Value2 = $ws.Range("A2").Value2
$ws.Range("A2").Value2 = $null

$ws.Range("A9").Value2 = $ws.Range("A6").Value2
$ws.Range("A6").Value2 = $null

$ws.Range("A15").Value2 = $ws.Range("A12").Value2
$ws.Range("A12").Value2 = $null

# Update view state: clear the frozen top-left cell and change the selection.
$ws.Range("B2").Select()
